# Add a new row (row 3) to Sheet1 with data for "Merge k Sorted Lists"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Populate the new row's values, matching the pattern of row 2
$ws.Range("A3").Value = "11/06/20201"
$ws.Range("B3").Value = 23
$ws.Range("C3").Value = "CHETHAN"
$ws.Range("D3").Value = "Merge k Sorted Lists"
$ws.Range("E3").Value = "Linked LIST"
$ws.Range("F3").Value = "HARD"

# Copy the formatting/style from row 2 so the new row matches (centered style)
$ws.Range("A2:F2").Copy() | Out-Null
$ws.Range("A3:F3").PasteSpecial(-4122) | Out-Null ; # xlPasteFormats
$excel.CutCopyMode = $false

# Update the selection to D3, matching the target state
$ws.Range("D3").Select() | Out-Null
